$wb = $excel.ActiveWorkbook

# --- Update the "StatOutput" sheet: stats for the Golden Retriever filter ---
# The cells originally hold numeric-looking text (shared strings), not real
# numbers, so force text entry (NumberFormat "@") then clear the format back
# off so no new cell style sticks around, while the value stays text-typed.
$statOutput = $wb.Worksheets.Item("StatOutput")

$statOutput.Range("A2").NumberFormat = "@"
$statOutput.Range("A2").Value = "29"
$statOutput.Range("A2").ClearFormats()

$statOutput.Range("B2").NumberFormat = "@"
$statOutput.Range("B2").Value = "25"
$statOutput.Range("B2").ClearFormats()

$statOutput.Range("C2").NumberFormat = "@"
$statOutput.Range("C2").Value = "24"
$statOutput.Range("C2").ClearFormats()

$statOutput.Range("D2").NumberFormat = "@"
$statOutput.Range("D2").Value = "2"
$statOutput.Range("D2").ClearFormats()

# --- Update the "StatOutput_Message" sheet: the Cypher query text (row 18) ---
$newQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Golden Retriever'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

$statOutputMessage = $wb.Worksheets.Item("StatOutput_Message")
$statOutputMessage.Range("A18").Value = $newQuery
